# Adds a "4-week low sales check" to the forecast summary workbook.
# - Recomputes MyForecast, Inventory Coverage, Stockout Risk, Reorder
#   Urgency and Seasonality Index on the "Forecast Comparison" sheet for
#   weeks 2-17 (only week 2 keeps a non-zero forecast / inventory figure).
# - Refreshes the dependent totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# ---------------------------------------------------------------------
# Forecast Comparison sheet (rows 2-17 correspond to weeks W10-W25)
# ---------------------------------------------------------------------

# Row 2 (W10) - still has inventory on hand / an urgent reorder
$wsForecast.Range("D2").Value = 4
$wsForecast.Range("H2").Value = 0.25
$wsForecast.Range("I2").Value = "High"
$wsForecast.Range("J2").Value = "Urgent"
$wsForecast.Range("L2").Value = 1.04

# Row 3 (W11)
$wsForecast.Range("D3").Value = 0
$wsForecast.Range("H3").Value = $null
$wsForecast.Range("I3").Value = "Low"
$wsForecast.Range("J3").Value = "Normal"
$wsForecast.Range("L3").Value = 1.13

# Row 4 (W12)
$wsForecast.Range("D4").Value = 0
$wsForecast.Range("H4").Value = $null
$wsForecast.Range("I4").Value = "Low"
$wsForecast.Range("J4").Value = "Normal"
$wsForecast.Range("L4").Value = 0.88

# Row 5 (W13)
$wsForecast.Range("D5").Value = 0
$wsForecast.Range("H5").Value = $null
$wsForecast.Range("I5").Value = "Low"
$wsForecast.Range("J5").Value = "Normal"
$wsForecast.Range("L5").Value = 1.08

# Row 6 (W14)
$wsForecast.Range("D6").Value = 0
$wsForecast.Range("H6").Value = $null
$wsForecast.Range("I6").Value = "Low"
$wsForecast.Range("J6").Value = "Normal"
$wsForecast.Range("L6").Value = 0.8100000000000001

# Row 7 (W15)
$wsForecast.Range("D7").Value = 0
$wsForecast.Range("H7").Value = $null
$wsForecast.Range("I7").Value = "Low"
$wsForecast.Range("J7").Value = "Normal"
$wsForecast.Range("L7").Value = 0.99

# Row 8 (W16)
$wsForecast.Range("D8").Value = 0
$wsForecast.Range("H8").Value = $null
$wsForecast.Range("I8").Value = "Low"
$wsForecast.Range("J8").Value = "Normal"
$wsForecast.Range("L8").Value = 0.98

# Row 9 (W17)
$wsForecast.Range("D9").Value = 0
$wsForecast.Range("H9").Value = $null
$wsForecast.Range("I9").Value = "Low"
$wsForecast.Range("J9").Value = "Normal"
$wsForecast.Range("L9").Value = 0.85

# Row 10 (W18)
$wsForecast.Range("D10").Value = 0
$wsForecast.Range("H10").Value = $null
$wsForecast.Range("I10").Value = "Low"
$wsForecast.Range("J10").Value = "Normal"
$wsForecast.Range("L10").Value = 1.1

# Row 11 (W19)
$wsForecast.Range("D11").Value = 0
$wsForecast.Range("H11").Value = $null
$wsForecast.Range("I11").Value = "Low"
$wsForecast.Range("J11").Value = "Normal"
$wsForecast.Range("L11").Value = 1.06

# Row 12 (W20)
$wsForecast.Range("D12").Value = 0
$wsForecast.Range("H12").Value = $null
$wsForecast.Range("I12").Value = "Low"
$wsForecast.Range("J12").Value = "Normal"
$wsForecast.Range("L12").Value = 1.15

# Row 13 (W21)
$wsForecast.Range("D13").Value = 0
$wsForecast.Range("H13").Value = $null
$wsForecast.Range("I13").Value = "Low"
$wsForecast.Range("J13").Value = "Normal"
$wsForecast.Range("L13").Value = 0.98

# Row 14 (W22)
$wsForecast.Range("D14").Value = 0
$wsForecast.Range("H14").Value = $null
$wsForecast.Range("I14").Value = "Low"
$wsForecast.Range("J14").Value = "Normal"
$wsForecast.Range("L14").Value = 0.92

# Row 15 (W23)
$wsForecast.Range("D15").Value = 0
$wsForecast.Range("H15").Value = $null
$wsForecast.Range("I15").Value = "Low"
$wsForecast.Range("J15").Value = "Normal"
$wsForecast.Range("L15").Value = 0.82

# Row 16 (W24)
$wsForecast.Range("D16").Value = 0
$wsForecast.Range("H16").Value = $null
$wsForecast.Range("I16").Value = "Low"
$wsForecast.Range("J16").Value = "Normal"
$wsForecast.Range("L16").Value = 1.2

# Row 17 (W25)
$wsForecast.Range("D17").Value = 0
$wsForecast.Range("H17").Value = $null
$wsForecast.Range("I17").Value = "Low"
$wsForecast.Range("J17").Value = "Normal"
$wsForecast.Range("L17").Value = 1.01

# ---------------------------------------------------------------------
# Summary sheet - totals derived from the new MyForecast column above
# (values are stored as text on this sheet, so force Text format before
# assigning, matching the original workbook layout).
# ---------------------------------------------------------------------

$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "4"

$wsSummary.Range("B10").NumberFormat = "@"
$wsSummary.Range("B10").Value = "4"

$wsSummary.Range("B11").NumberFormat = "@"
$wsSummary.Range("B11").Value = "4"

$wsSummary.Range("B12").NumberFormat = "@"
$wsSummary.Range("B12").Value = "4"

$wsSummary.Range("B14").NumberFormat = "@"
$wsSummary.Range("B14").Value = "0"
